# Update the "魔法点获取的比率" (EnergyRate) column (F) on the "Job" sheet
# so the card-region ratios show correctly in the deck view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Job")

$ws.Range("F4").Value2  = "75;10;15"
$ws.Range("F5").Value2  = "60;20;20"
$ws.Range("F6").Value2  = "60;20;20"
$ws.Range("F7").Value2  = "60;15;25"
$ws.Range("F8").Value2  = "60;25;15"
$ws.Range("F9").Value2  = "55;20;25"
$ws.Range("F10").Value2 = "50;10;40"
$ws.Range("F11").Value2 = "50;15;35"
$ws.Range("F12").Value2 = "55;10;35"
$ws.Range("F13").Value2 = "55;10;35"
$ws.Range("F14").Value2 = "50;20;30"

# Leave F15/F16 (活动1/活动2) untouched - values unchanged in the source edit.

# Match the last-selected cell recorded in the saved workbook.
[void]$ws.Range("F9").Select()
